# Updated cryptos list on Tue Aug  1 08:53:42 UTC 2023 with GitHub Actions
#
# This script refreshes the Price (column D) and Volume(1h) (column E)
# values for the rows in the crypto table. Every value is kept as plain
# text (matching the original inlineStr cells), so numeric-looking
# strings (e.g. "0.9996", "1.0000") must NOT be auto-converted into
# numbers by Excel. We do that by writing the text into a scratch cell
# via a formula that evaluates to a text string (so the cell's data type
# is "text" without needing any NumberFormat/quote-prefix styling), then
# copy/paste-special just the *values* onto the destination cell. This
# changes the destination's stored value/type without touching its style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $targetAddr, $text)
    $scratch = $ws.Range("Z100")
    $escaped = $text.Replace("""", """""")
    $scratch.Formula = "=""" + $escaped + """"
    $scratch.Copy()
    $dst = $ws.Range($targetAddr)
    $dst.PasteSpecial(-4163)   # xlPasteValues: values only, no formatting
    $scratch.ClearContents()
}

Set-TextValue $ws "D2"  '28.945.55'
Set-TextValue $ws "E2"  '  -1.58%  '
Set-TextValue $ws "D3"  '1.833.84'
Set-TextValue $ws "D4"  '0.9996'
Set-TextValue $ws "E4"  '  -0.15%  '
Set-TextValue $ws "D5"  '245.86'
Set-TextValue $ws "E5"  '  +0.77%  '
Set-TextValue $ws "D6"  '0.6897'
Set-TextValue $ws "E6"  '  -2.18%  '
Set-TextValue $ws "D7"  '1.0000'
Set-TextValue $ws "E7"  '  -0.11%  '
Set-TextValue $ws "D8"  '0.07691'
Set-TextValue $ws "E8"  '  -2.93%  '
Set-TextValue $ws "D9"  '0.3052'
Set-TextValue $ws "E9"  '  -2.80%  '
Set-TextValue $ws "E10" '  -4.08%  '
Set-TextValue $ws "D11" '0.07799'
Set-TextValue $ws "E11" '  -1.23%  '
Set-TextValue $ws "D12" '1.841.22'
Set-TextValue $ws "E12" '  -1.44%  '
Set-TextValue $ws "D13" '5.079'
Set-TextValue $ws "E13" '  -2.19%  '
Set-TextValue $ws "D14" '90.49'
Set-TextValue $ws "E14" '  -3.62%  '
Set-TextValue $ws "D15" '0.6801'
Set-TextValue $ws "E15" '  -3.24%  '
Set-TextValue $ws "D16" '6.436'
Set-TextValue $ws "E16" '  -1.49%  '
Set-TextValue $ws "D17" '0.000008338'
Set-TextValue $ws "E17" '  -0.68%  '
Set-TextValue $ws "D18" '28.957.25'
Set-TextValue $ws "E18" '  -1.52%  '
Set-TextValue $ws "D19" '243.37'
Set-TextValue $ws "E19" '  -4.24%  '
Set-TextValue $ws "D20" '2.084.40'
Set-TextValue $ws "E20" '  -1.93%  '
Set-TextValue $ws "D21" '12.72'
Set-TextValue $ws "E21" '  -3.06%  '
Set-TextValue $ws "D22" '0.9996'
Set-TextValue $ws "D23" '7.473'
Set-TextValue $ws "E23" '  -2.30%  '
Set-TextValue $ws "D24" '1.000'
Set-TextValue $ws "E24" '  -0.09%  '
Set-TextValue $ws "D25" '162.21'
Set-TextValue $ws "E25" '  +0.60%  '
Set-TextValue $ws "D26" '0.1469'
Set-TextValue $ws "E26" '  -5.84%  '
Set-TextValue $ws "E27" '  -2.36%  '
Set-TextValue $ws "D28" '18.21'
Set-TextValue $ws "E28" '  -3.21%  '
Set-TextValue $ws "E29" '  +3.11%  '
Set-TextValue $ws "D30" '4.211'
Set-TextValue $ws "E30" '  -2.75%  '
Set-TextValue $ws "D31" '4.155'
Set-TextValue $ws "E31" '  -2.55%  '
Set-TextValue $ws "D32" '1.177'
Set-TextValue $ws "E32" '  -3.10%  '
Set-TextValue $ws "D33" '0.05124'
Set-TextValue $ws "E33" '  -3.31%  '
Set-TextValue $ws "D34" '0.7660'
Set-TextValue $ws "E34" '  +1.91%  '
Set-TextValue $ws "D35" '1.846'
Set-TextValue $ws "E35" '  -2.54%  '
Set-TextValue $ws "D36" '1.144'
Set-TextValue $ws "E36" '  -2.73%  '
Set-TextValue $ws "E37" '  -1.17%  '
Set-TextValue $ws "D39" '1.234.41'
Set-TextValue $ws "E39" '  -4.15%  '
Set-TextValue $ws "D40" '2.696'
Set-TextValue $ws "E40" '  -2.57%  '
Set-TextValue $ws "D41" '0.9247'
Set-TextValue $ws "E41" '  +3.35%  '
Set-TextValue $ws "D42" '108.58'
Set-TextValue $ws "E42" '  -0.59%  '
Set-TextValue $ws "D43" '5.834'
Set-TextValue $ws "E43" '  -3.21%  '
Set-TextValue $ws "D44" '0.9994'
Set-TextValue $ws "D45" '9.589'
Set-TextValue $ws "E45" '  +0.03%  '
Set-TextValue $ws "E46" '  -3.46%  '
Set-TextValue $ws "D47" '1.982.75'
Set-TextValue $ws "E47" '  -2.18%  '
Set-TextValue $ws "E48" '  -0.16%  '
Set-TextValue $ws "D49" '64.21'
Set-TextValue $ws "E49" '  -9.80%  '
Set-TextValue $ws "E50" '  -2.94%  '
Set-TextValue $ws "D51" '6.930'
Set-TextValue $ws "E51" '  -2.03%  '

$excel.CutCopyMode = $false
